# Applies the "update presentation about the project" edits.
#
# Strategy: for every text change we locate the exact original run text
# (the full text of a single <a:r>/<a:t> run) inside the shape's
# TextRange via IndexOf, then rewrite that exact character span with
# Characters(start, length).Text = "...". Using the *entire* original
# run's text as the search/replace span (instead of just the changed
# fragment) keeps the run count/formatting identical to the source -
# only the <a:t> contents change, exactly like the diff.

$p = $ppt.ActivePresentation

function Set-RunText {
    param(
        $TextRange,
        [string]$OldText,
        [string]$NewText
    )
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Could not find text: [$OldText]"
    }
    $chars = $TextRange.Characters($idx + 1, $OldText.Length)
    $chars.Text = $NewText
}

# ---------------------------------------------------------------------
# Slide 2 - "La coda" content placeholder (Shapes.Item(2))
# ---------------------------------------------------------------------
# NOTE: TextRange.Text normalizes the typographic right-single-quote
# (U+2019) down to a plain ASCII apostrophe when *read back*, so the
# search ("old") strings below use a plain apostrophe; the replacement
# ("new") strings use the real U+2019 so the saved XML matches the diff.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange

Set-RunText $tr2 "Questi nodi vengono gestiti all'interno dell'oggetto " "Questi nodi vengono gestiti all’interno della classe "
Set-RunText $tr2 " che esegue diverse funzioni sui nodi (in maniera non " " che esegue diverse funzioni sui nodi (in modo non "
Set-RunText $tr2 "[reading] -> nodo in cima alla coda, senza eliminarlo" "[reading] -> primo nodo della coda (come pop), ma senza eliminarlo"

# ---------------------------------------------------------------------
# Slide 3 - "Policy" content placeholder (Shapes.Item(2))
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange

# Resize/reposition the placeholder and drop the auto-shrink-text-on-overflow.
$sh3.Top = 106.4582748413086
$sh3.Height = 327.0834045410156
$sh3.TextFrame.AutoSize = 0

Set-RunText $tr3 "L'oggetto Policy generalizza una politica di scheduling rendendo i metodi della coda " "La classe Policy generalizza una politica di scheduling rendendo i metodi della coda "
Set-RunText $tr3 " (importati dalla libreria threading)." " (importati dalla libreria threading). Con i vari oggetti figli di Policy è poi possibile invocare i metodi visti in precedenza per interagire con la coda."
Set-RunText $tr3 " in base alla politica adottata." " in base alla politica adottata. Tali Funzioni sono state implementate negli oggetti Sem delle varie politiche."

# ---------------------------------------------------------------------
# Slide 4 - title textbox "Implementazione Policy" -> "Tipi di Policy"
# (no position change on this slide)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(6).TextFrame.TextRange
Set-RunText $tr4 "Implementazione Policy" "Tipi di Policy"

# ---------------------------------------------------------------------
# Slide 5 - title textbox, reposition + rename
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(6)
$sh5.Top = 16.675907135009766
$tr5 = $sh5.TextFrame.TextRange
Set-RunText $tr5 "Implementazione Policy" "Tipi di Policy"

# ---------------------------------------------------------------------
# Slide 6 - title textbox, reposition + rename
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(6)
$sh6.Top = 15.704252243041992
$tr6 = $sh6.TextFrame.TextRange
Set-RunText $tr6 "Implementazione Policy" "Tipi di Policy"
